$d = $word.ActiveDocument

# Insert the new "Mascarillas y protectores faciales" section (heading,
# three bulleted paragraphs and a 6x3 data table) right before the final
# paragraph mark of the document body, i.e. immediately after the last
# existing paragraph and before the sectPr.
$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$sectionXml = @'
<w:p><w:pPr><w:pStyle w:val="Ttulo1"/></w:pPr><w:r><w:t>2. Mascarillas y protectores faciales</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaconvietas"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Mediante el Decreto de Urgencia N° 021-2021 y la Resolución de Secretaría General N° 047-2021-MINEDU, se transfirieron S/ </w:t></w:r><w:r><w:t>1.1</w:t></w:r><w:r><w:t xml:space="preserve"> millones de soles para la adquisición y distribución de mascarillas faciales textiles de uso comunitario para estudiantes y personal que labora en instituciones educativas públicas, así como protectores faciales para el mencionado personal.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaconvietas"/><w:jc w:val="both"/></w:pPr><w:r><w:t>La adquisición de mascarillas y protectores faciales es condición necesaria para el retorno seguro a los servicios educativos presenciales y semipresenciales, según lo dispuesto por las “Disposiciones para la prestación del servicio en las instituciones y programas educativos públicos y privados de la Educación Básica de los ámbitos urbanos y rurales, en el marco de la emergencia sanitaria de la COVID-19”, aprobado mediante Resolución Ministerial N° 121-2021- MINEDU y modificado con Resoluciones Ministeriales N° 199-2021-MINEDU y N° 273-2021- MINEDU.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaconvietas"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Con fecha de corte al </w:t></w:r><w:r><w:t>21 de setiembre de 2021</w:t></w:r><w:r><w:t xml:space="preserve">, la ejecución a nivel regional de los recursos de mascarillas faciales textiles protectores faciales fue del </w:t></w:r><w:r><w:t>92.1%</w:t></w:r><w:r><w:t xml:space="preserve"> (devengado) según se presenta a continuación:</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="Listavistosa-nfasis1"/><w:tblW w:type="auto" w:w="0"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="1440"/><w:gridCol w:w="1440"/><w:gridCol w:w="1440"/><w:gridCol w:w="1440"/><w:gridCol w:w="1440"/><w:gridCol w:w="1440"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>region</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>nom_ue</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>certificado</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>comprometido_anual</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>devengado</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>transferencia</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>TACNA</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>300. EDUCACION TACNA</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>79465.3984375</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>79465.3984375</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>79465.3984375</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>120391</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>TACNA</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>301. UGEL TACNA</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>898014.8125</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>898014.8125</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>898014.8125</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="1440"/></w:tcPr><w:p><w:r><w:t>940771</w:t></w:r></w:p></w:tc></w:tr></w:tbl>
'@

$insertionPoint.InsertXML($sectionXml)
